$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark from the UberLibs bullet (it sits
#    collapsed right between "...create wishlists" and ", borrow books...").
#    Bookmarks.Delete() is not reliably wired up in this host, but deleting a
#    tiny range that straddles the bookmark's position and re-typing the same
#    text does remove it, confirmed by Bookmarks.Exists() afterwards.
# ---------------------------------------------------------------------------
$locator = $d.Content
$found = $locator.Find.Execute("create wishlists, borrow books", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
if ($found -and $d.Bookmarks.Exists("_GoBack")) {
    $bmPos = $locator.Start + "create wishlists".Length
    $span = $d.Range($bmPos - 1, $bmPos + 1)
    $spanText = $span.Text
    $span.Delete()
    $collapsed = $d.Range($bmPos - 1, $bmPos - 1)
    $collapsed.InsertAfter($spanText)
}

# ---------------------------------------------------------------------------
# 2. Update the LG Soft India bullet: "LG's next generation JavaScript
#    framework" -> "LG's next generation React based JS framework".
#    Locate the precise phrase first (so we don't touch any of the other
#    "JavaScript" mentions elsewhere in the resume), then replace just the
#    "JavaScript" word via a forward-only, non-wrapping Find scoped to start
#    right after "LG's next generation " - this keeps the surrounding runs
#    ("Researched ..." etc.) untouched instead of retyping the whole match.
# ---------------------------------------------------------------------------
$locator2 = $d.Content
$found2 = $locator2.Find.Execute("LG's next generation JavaScript framework", $true, $false, $false, $false, `
                                  $false, $true, 1, $false, "", 0)
if ($found2) {
    $wordStart = $locator2.Start + "LG's next generation ".Length
    $scoped = $d.Range($wordStart, $d.Content.End)
    $found3 = $scoped.Find.Execute("JavaScript", $true, $true, $false, $false, $false, `
                                    $true, 0, $false, "React based JS", 1)

    # -----------------------------------------------------------------------
    # 3. Re-create the "_GoBack" bookmark right after the newly-typed text,
    #    matching where Word leaves it after the most recent edit.
    # -----------------------------------------------------------------------
    if ($found3) {
        $newBmPos = $wordStart + "React based JS".Length
        $newBmRange = $d.Range($newBmPos, $newBmPos)
        $d.Bookmarks.Add("_GoBack", $newBmRange)
    }
}
